# Applies updated market-price / profit figures to the Leve profit tables
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per the
# scheduled-runner data refresh described in the commit.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 45494.668
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 45494.668
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 45494.668
$ws.Range("N3").Value = -45722.668
$ws.Range("H12").Value = 125
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -490
$ws.Range("H15").Value = 6150.224
$ws.Range("I15").Value = 6150.224
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 18450.672
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -18281.672
$ws.Range("H33").Value = 258.16278
$ws.Range("I33").Value = 151.91891
$ws.Range("J33").Value = 913.3333
$ws.Range("K33").Value = 151.91891
$ws.Range("L33").Value = 913.3333
$ws.Range("M33").Value = 77.08108999999999
$ws.Range("N33").Value = -1371.3333
$ws.Range("H34").Value = 22619.9
$ws.Range("I34").Value = 17733.334
$ws.Range("J34").Value = 29949.75
$ws.Range("K34").Value = 17733.334
$ws.Range("L34").Value = 29949.75
$ws.Range("M34").Value = -17530.334
$ws.Range("N34").Value = -30355.75
$ws.Range("H36").Value = 22619.9
$ws.Range("I36").Value = 17733.334
$ws.Range("J36").Value = 29949.75
$ws.Range("K36").Value = 17733.334
$ws.Range("L36").Value = 29949.75
$ws.Range("M36").Value = -17018.334
$ws.Range("N36").Value = -31379.75
$ws.Range("H40").Value = 7928.9414
$ws.Range("I40").Value = 10283.5
$ws.Range("J40").Value = 2278
$ws.Range("K40").Value = 10283.5
$ws.Range("L40").Value = 2278
$ws.Range("M40").Value = -10108.5
$ws.Range("N40").Value = -2628
$ws.Range("H76").Value = 3087.1875
$ws.Range("I76").Value = 2999.2222
$ws.Range("J76").Value = 3200.2856
$ws.Range("K76").Value = 2999.2222
$ws.Range("L76").Value = 3200.2856
$ws.Range("M76").Value = -2684.2222
$ws.Range("N76").Value = -3830.2856
$ws.Range("H79").Value = 3087.1875
$ws.Range("I79").Value = 2999.2222
$ws.Range("J79").Value = 3200.2856
$ws.Range("K79").Value = 2999.2222
$ws.Range("L79").Value = 3200.2856
$ws.Range("M79").Value = -1907.2222
$ws.Range("N79").Value = -5384.2856
$ws.Range("H92").Value = 322.73077
$ws.Range("I92").Value = 256.17392
$ws.Range("J92").Value = 833
$ws.Range("K92").Value = 256.17392
$ws.Range("L92").Value = 833
$ws.Range("M92").Value = 991.82608
$ws.Range("N92").Value = -3329
$ws.Range("H94").Value = 2591.3635
$ws.Range("I94").Value = 2833.889
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 2833.889
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -2382.889
$ws.Range("N94").Value = -2402
$ws.Range("H100").Value = 1494.0588
$ws.Range("I100").Value = 870
$ws.Range("J100").Value = 2385.5715
$ws.Range("K100").Value = 870
$ws.Range("L100").Value = 2385.5715
$ws.Range("M100").Value = -329
$ws.Range("N100").Value = -3467.5715
$ws.Range("H101").Value = 1167.25
$ws.Range("I101").Value = 695.6
$ws.Range("J101").Value = 1953.3334
$ws.Range("K101").Value = 2086.8
$ws.Range("L101").Value = 5860.0002
$ws.Range("M101").Value = -464.8000000000002
$ws.Range("N101").Value = -9104.0002
$ws.Range("H102").Value = 45494.668
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 45494.668
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 45494.668
$ws.Range("N102").Value = -51984.668
$ws.Range("H103").Value = 1306.6666
$ws.Range("I103").Value = 1028.75
$ws.Range("J103").Value = 1862.5
$ws.Range("K103").Value = 3086.25
$ws.Range("L103").Value = 5587.5
$ws.Range("M103").Value = -2500.25
$ws.Range("N103").Value = -6759.5
$ws.Range("H105").Value = 48823.668
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 48823.668
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 48823.668
$ws.Range("N105").Value = -55811.668
$ws.Range("H132").Value = 23126.455
$ws.Range("I132").Value = 3583.0303
$ws.Range("J132").Value = 81756.73
$ws.Range("K132").Value = 10749.0909
$ws.Range("L132").Value = 245270.19
$ws.Range("M132").Value = -8219.090899999999
$ws.Range("N132").Value = -250330.19
$ws.Range("H135").Value = 12501030
$ws.Range("I135").Value = 682.14703
$ws.Range("J135").Value = 83336340
$ws.Range("K135").Value = 6139.32327
$ws.Range("L135").Value = 750027060
$ws.Range("M135").Value = -3604.32327
$ws.Range("N135").Value = -750032130
$ws.Range("H141").Value = 2327
$ws.Range("I141").Value = 854.6667
$ws.Range("J141").Value = 8952.5
$ws.Range("K141").Value = 2564.0001
$ws.Range("L141").Value = 26857.5
$ws.Range("M141").Value = 2615.9999
$ws.Range("N141").Value = -37217.5

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1300.5555
$ws.Range("I97").Value = 1440.6666
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 1440.6666
$ws.Range("L97").Value = 600
$ws.Range("M97").Value = -944.6666
$ws.Range("N97").Value = -1592
$ws.Range("H102").Value = 41988.4
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 52235.5
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 52235.5
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -55479.5
$ws.Range("H106").Value = 45925
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 45925
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 45925
$ws.Range("N106").Value = -48449

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 40492
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 40492
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 40492
$ws.Range("N95").Value = -45984
$ws.Range("H99").Value = 1920.2051
$ws.Range("I99").Value = 1861
$ws.Range("J99").Value = 2070.9092
$ws.Range("K99").Value = 1861
$ws.Range("L99").Value = 2070.9092
$ws.Range("M99").Value = -363
$ws.Range("N99").Value = -5066.9092
$ws.Range("H100").Value = 42653.332
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 42653.332
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 42653.332
$ws.Range("N100").Value = -44817.332
$ws.Range("H105").Value = 2527.9656
$ws.Range("I105").Value = 1817.5
$ws.Range("J105").Value = 4106.778
$ws.Range("K105").Value = 1817.5
$ws.Range("L105").Value = 4106.778
$ws.Range("M105").Value = -70.5
$ws.Range("N105").Value = -7600.778

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1466.6666
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -1213
$ws.Range("N16").Value = -1974
$ws.Range("H22").Value = 3798.5
$ws.Range("I22").Value = 495.5
$ws.Range("J22").Value = 5450
$ws.Range("K22").Value = 495.5
$ws.Range("L22").Value = 5450
$ws.Range("M22").Value = -145.5
$ws.Range("N22").Value = -6150
$ws.Range("H86").Value = 3091.182
$ws.Range("I86").Value = 3613.4285
$ws.Range("J86").Value = 2177.25
$ws.Range("K86").Value = 3613.4285
$ws.Range("L86").Value = 2177.25
$ws.Range("M86").Value = -2490.4285
$ws.Range("N86").Value = -4423.25
$ws.Range("H89").Value = 3091.182
$ws.Range("I89").Value = 3613.4285
$ws.Range("J89").Value = 2177.25
$ws.Range("K89").Value = 18067.1425
$ws.Range("L89").Value = 10886.25
$ws.Range("M89").Value = -12451.1425
$ws.Range("N89").Value = -22118.25
$ws.Range("H92").Value = 40424.855
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 40424.855
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 40424.855
$ws.Range("N92").Value = -45416.855
$ws.Range("H113").Value = 1466.6666
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -5740

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 446.7619
$ws.Range("I23").Value = 523.6667
$ws.Range("J23").Value = 389.08334
$ws.Range("K23").Value = 1571.0001
$ws.Range("L23").Value = 1167.25002
$ws.Range("M23").Value = -1336.0001
$ws.Range("N23").Value = -1637.25002
$ws.Range("H49").Value = 1600
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1600
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 4800
$ws.Range("N49").Value = -5112
$ws.Range("H56").Value = 122708.586
$ws.Range("I56").Value = 122708.586
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 122708.586
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -122178.586
$ws.Range("H97").Value = 1346.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1346.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4039.5
$ws.Range("N97").Value = -5031.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4284.645
$ws.Range("I70").Value = 4268.96
$ws.Range("J70").Value = 4350
$ws.Range("K70").Value = 4268.96
$ws.Range("L70").Value = 4350
$ws.Range("M70").Value = -3998.96
$ws.Range("N70").Value = -4890
$ws.Range("H73").Value = 4284.645
$ws.Range("I73").Value = 4268.96
$ws.Range("J73").Value = 4350
$ws.Range("K73").Value = 4268.96
$ws.Range("L73").Value = 4350
$ws.Range("M73").Value = -3332.96
$ws.Range("N73").Value = -6222
$ws.Range("H97").Value = 3726.5417
$ws.Range("I97").Value = 3031
$ws.Range("J97").Value = 4885.778
$ws.Range("K97").Value = 3031
$ws.Range("L97").Value = 4885.778
$ws.Range("M97").Value = -2535
$ws.Range("N97").Value = -5877.778
$ws.Range("H116").Value = 38998
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 38998
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 38998
$ws.Range("N116").Value = -48176

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 48063.25
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 48063.25
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 48063.25
$ws.Range("N36").Value = -49187.25
$ws.Range("H46").Value = 6195.5713
$ws.Range("I46").Value = 1251.1428
$ws.Range("J46").Value = 11140
$ws.Range("K46").Value = 1251.1428
$ws.Range("L46").Value = 11140
$ws.Range("M46").Value = -1063.1428
$ws.Range("N46").Value = -11516
$ws.Range("H93").Value = 1493.4231
$ws.Range("I93").Value = 992.7857
$ws.Range("J93").Value = 2077.5
$ws.Range("K93").Value = 992.7857
$ws.Range("L93").Value = 2077.5
$ws.Range("M93").Value = 255.2143
$ws.Range("N93").Value = -4573.5
$ws.Range("H100").Value = 2791.7693
$ws.Range("I100").Value = 2789.3
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 2789.3
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -2248.3
$ws.Range("N100").Value = -3882
$ws.Range("H104").Value = 36340
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 36340
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 36340
$ws.Range("N104").Value = -43328
$ws.Range("H105").Value = 50292
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 50292
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 50292
$ws.Range("N105").Value = -57280

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 473
$ws.Range("I100").Value = 297.875
$ws.Range("J100").Value = 940
$ws.Range("K100").Value = 595.75
$ws.Range("L100").Value = 1880
$ws.Range("M100").Value = -54.75
$ws.Range("N100").Value = -2962
$ws.Range("H101").Value = 21850.2
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 21850.2
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 21850.2
$ws.Range("N101").Value = -28340.2
$ws.Range("H120").Value = 47416
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 47416
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 47416
$ws.Range("N120").Value = -57092

Write-Host "Updated leve profit figures across all sheets."
